$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'300.35"
$ws.Range("E2").Value = "'2.68%"
$ws.Range("D3").Value = "'42.31"
$ws.Range("E3").Value = "'4.46%"
$ws.Range("D4").Value = "'5.027"
$ws.Range("E4").Value = "'0.29%"
$ws.Range("D5").Value = "'0.07617"
$ws.Range("E5").Value = "'3.01%"
$ws.Range("D6").Value = "'1.606"
$ws.Range("E6").Value = "'1.76%"
$ws.Range("D7").Value = "'0.9763"
$ws.Range("E7").Value = "'5.72%"
$ws.Range("D9").Value = "'0.1200"
$ws.Range("E9").Value = "'0.90%"
$ws.Range("D10").Value = "'0.1832"
$ws.Range("E10").Value = "'1.02%"
$ws.Range("D11").Value = "'0.09125"
$ws.Range("E11").Value = "'3.99%"
$ws.Range("E12").Value = "'-4.69%"
$ws.Range("E13").Value = "'-0.57%"
$ws.Range("D14").Value = "'0.001262"
$ws.Range("E14").Value = "'-0.59%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005929"
$ws.Range("E15").Value = "'1.60%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.354"
$ws.Range("E16").Value = "'0.36%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.382"
$ws.Range("E17").Value = "'2.06%"
$ws.Range("B18").Value = "BitpandaEcosystemToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D18").Value = "'0.3341"
$ws.Range("E18").Value = "'0.71%"
$ws.Range("B19").Value = "MCDex"
$ws.Range("C19").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D19").Value = "'8.343"
$ws.Range("E19").Value = "'5.43%"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "'0.1347"
$ws.Range("E20").Value = "'-3.14%"
$ws.Range("B21").Value = "ZBToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D21").Value = "'0.3290"
$ws.Range("E21").Value = "'11.16%"
$ws.Range("B22").Value = "CoinExToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D22").Value = "'0.04113"
$ws.Range("E22").Value = "'4.97%"
$ws.Range("D23").Value = "'0.001264"
$ws.Range("E23").Value = "'0.19%"
$ws.Range("D24").Value = "'0.003938"
$ws.Range("E24").Value = "'3.49%"
$ws.Range("D25").Value = "'0.0001345"
$ws.Range("E25").Value = "'9.24%"
$ws.Range("D38").Value = "'0.02424"
$ws.Range("E38").Value = "'4.19%"
$ws.Range("D39").Value = "'0.05256"
$ws.Range("E39").Value = "'3.80%"
$ws.Range("E40").Value = "'11.05%"
$ws.Range("D41").Value = "'0.007684"
$ws.Range("E41").Value = "'-1.74%"
$ws.Range("D42").Value = "'0.1341"
$ws.Range("E42").Value = "'4.02%"
$ws.Range("D43").Value = "'0.007372"
$ws.Range("E43").Value = "'-0.25%"
$ws.Range("D44").Value = "'0.007277"
$ws.Range("E44").Value = "'-9.67%"
$ws.Range("D45").Value = "'0.3008"
$ws.Range("E45").Value = "'3.32%"
$ws.Range("D46").Value = "'0.00006334"
$ws.Range("E46").Value = "'1.76%"
$ws.Range("D47").Value = "'0.00000000748"
$ws.Range("E47").Value = "'-0.39%"
$ws.Range("D48").Value = "'0.04604"
$ws.Range("E48").Value = "'-3.46%"
$ws.Range("E49").Value = "'-0.06%"
$ws.Range("D50").Value = "'0.00002093"
$ws.Range("E50").Value = "'-0.39%"
$ws.Range("D51").Value = "'0.0001994"
$ws.Range("E51").Value = "'-0.39%"
